$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 0.1.0
$wsMeta.Range("B3").Value = "0.1.0"

# Status: active -> draft
$wsMeta.Range("B6").Value = "draft"

# Experimental: (blank) -> false
$wsMeta.Range("B7").Value = "false"

# Date: 2025-11-28T01:24:36+00:00 -> 2025-12-26T14:13:58+00:00
$wsMeta.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: (blank) -> Value set for categorizing nursing problems
$wsMeta.Range("B11").Value = "Value set for categorizing nursing problems"
